$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers; force Text format so
# Excel stores them as text (matching the source inlineStr cells) instead
# of silently converting them to numbers (which would also lose things
# like trailing zeros, e.g. "1.70" -> 1.7).
$numericLookingCells = @(
    "D5", "D6", "D7", "D10", "D11", "D14", "D19", "D21",
    "D22", "D24", "D27", "D28", "D30", "D31", "D32", "D36",
    "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46",
    "D47", "D49", "D50", "D51"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated crypto price / link / volume figures.
$ws.Range("D2").Value = '43.853.88'
$ws.Range("E2").Value = '  +2.98%  '
$ws.Range("D3").Value = '2.337.88'
$ws.Range("E3").Value = '  +2.22%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '311.99'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Value = '108.28'
$ws.Range("E6").Value = '  +2.98%  '
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").Value = '  +1.88%  '
$ws.Range("D10").Value = '41.23'
$ws.Range("E10").Value = '  +3.89%  '
$ws.Range("D11").Value = '0.0918'
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").Value = '1.01'
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("D16").Value = '2.694.24'
$ws.Range("E16").Value = '  +2.08%  '
$ws.Range("D17").Value = '2.333.44'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '43.778.23'
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = '7.54'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '13.04'
$ws.Range("E21").Value = '  -5.96%  '
$ws.Range("D22").Value = '74.14'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  -3.47%  '
$ws.Range("D24").Value = '268.71'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("E25").Value = '  +2.00%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '7.66'
$ws.Range("E27").Value = '  +6.06%  '
$ws.Range("D28").Value = '11.12'
$ws.Range("E28").Value = '  +2.64%  '
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("D30").Value = '39.09'
$ws.Range("E30").Value = '  +5.15%  '
$ws.Range("D31").Value = '22.58'
$ws.Range("D32").Value = '168.75'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("E34").Value = '  +9.32%  '
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").Value = '0.115'
$ws.Range("E36").Value = '  +1.23%  '
$ws.Range("E37").Value = '  +3.88%  '
$ws.Range("D38").Value = '0.0362'
$ws.Range("E38").Value = '  +3.10%  '
$ws.Range("D39").Value = '2.88'
$ws.Range("E39").Value = '  +8.57%  '
$ws.Range("D40").Value = '3.78'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").Value = '1.70'
$ws.Range("E41").Value = '  +8.24%  '
$ws.Range("D42").Value = '105.06'
$ws.Range("E42").Value = '  +10.73%  '
$ws.Range("E43").Value = '  +2.37%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D44").Value = '13.41'
$ws.Range("E44").Value = '  +10.33%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").Value = '71.47'
$ws.Range("E45").Value = '  +0.74%  '
$ws.Range("D46").Value = '1.01'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '114.23'
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").Value = '1.672.87'
$ws.Range("E48").Value = '  -3.97%  '
$ws.Range("D49").Value = '76.72'
$ws.Range("E49").Value = '  -4.08%  '
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = '0.217'
$ws.Range("E50").Value = '  +14.40%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '8.93'
$ws.Range("E51").Value = '  +2.21%  '
